# ex9.1.9(Linear)__M_Stationarygenerator_alpha_zero.xlsx
# "nuevos experimentos no convexos"
#
# The workbook stores every cell (even the ones that look numeric) as a
# plain text/shared-string value, so every write below forces the cell to
# Text format before assigning the value (otherwise Excel would silently
# re-interpret e.g. "2.215" as the number 2.215) and then restores the
# cell's style back to Normal so no stray number-format is left behind.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ---- Restricciones_del_follower ----------------------------------------
$ws = $wb.Worksheets.Item("Restricciones_del_follower")

Set-TextValue $ws.Range("A2") "7.865 - x - 0.5y"
Set-TextValue $ws.Range("B2") "-5.865"
Set-TextValue $ws.Range("D2") "0.93"
Set-TextValue $ws.Range("F2") "3.4000000000000004"

Set-TextValue $ws.Range("A3") "-4.975 - 0.25x + y"
Set-TextValue $ws.Range("B3") "2.9749999999999996"
Set-TextValue $ws.Range("D3") "0.41"
Set-TextValue $ws.Range("E3") "5.0"
Set-TextValue $ws.Range("F3") "0"

Set-TextValue $ws.Range("A4") "-7.865 + x + 0.5y"
Set-TextValue $ws.Range("B4") "-0.1349999999999998"
Set-TextValue $ws.Range("D4") "0.7"
Set-TextValue $ws.Range("F4") "6.6000000000000005"

Set-TextValue $ws.Range("A5") "-11.559999999999999 + x - 2y"
Set-TextValue $ws.Range("B5") "-9.559999999999999"
Set-TextValue $ws.Range("D5") "0.36"
Set-TextValue $ws.Range("E5") "7.9"
Set-TextValue $ws.Range("F5") "0"

Set-TextValue $ws.Range("A6") "-6.17 - y"
Set-TextValue $ws.Range("B6") "-6.17"
Set-TextValue $ws.Range("D6") "0.79"
Set-TextValue $ws.Range("E6") "0"
Set-TextValue $ws.Range("F6") "2.9"

# ---- Punto_modificado ----------------------------------------------------
$ws = $wb.Worksheets.Item("Punto_modificado")
Set-TextValue $ws.Range("A2") "4.78"
Set-TextValue $ws.Range("B2") "6.17"

# ---- Vector_bf ------------------------------------------------------------
$ws = $wb.Worksheets.Item("Vector_bf")
Set-TextValue $ws.Range("A2") "2.215"

# ---- Vector_BF ------------------------------------------------------------
# NB: sheet-name lookup is case-insensitive, and this workbook also has a
# sheet literally named "Vector_bf" - Worksheets.Item("Vector_BF") would
# resolve to that *other* sheet instead. Use the (1-based) tab index, which
# is unambiguous, to reach the real "Vector_BF" sheet (tab #6).
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-7.65"
Set-TextValue $ws.Range("A3") "9.8"
